# Update SegmentsTestData: modify test data values on the CreateSegment sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CreateSegment")

$ws.Range("D2").Value = "MT AUTOMATION"
$ws.Range("F2").Value = "mt_dc_esp"
$ws.Range("I2").Value = "Google Cookie;MAID"
$ws.Range("L2").Value = "Test Destination"
